# Diane's correction prior to discussion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - "Final field test" definition
$ws.Range("B3").Value = "Visits to households for practice interviews at the end of fieldworker training, to assess the trainees' readiness for fieldwork, which often also serves as a final test of all survey instruments and protocols."

# Row 4 - "Questionnaire" definition
$ws.Range("B4").Value = "A questionnaire is a research instrument consisting of a series of questions for the purpose of gathering information from respondents."

# Row 6 - "Data monitor" definition
$ws.Range("B6").Value = "The role of those conducting quality assurance checks of incoming fieldwork data."

# Row 7 - "Field practice" definition (also grows the wrapped row height)
$ws.Range("B7").Value = "Visits to households during fieldwork training, so that trainees can practice interview techniques."
$ws.Rows.Item(7).RowHeight = 30.75

# Row 10 - "Pre-test" definition
$ws.Range("B10").Value = "Evaluates the software (CAPI or the data entry program, if using PAPI) to ensure that survey data are entered correctly, the questionnaire flow is consistent, and data checks are complete. This can be done with a very small sample – even just a handful of households – as long as all the modules are administered."

# Row 12 - "Data collection mode" definition
$ws.Range("B12").Value = "The way survey data are collected. Traditionally, most surveys were conducted using PAPI, Paper Assisted Personal Interviewing, where responses were recorded on paper questionnaires and entered centrally after fieldwork. In, CAFE, Computer Assisted Field Entry, data entry happens in the field during fieldwork, usually to allow timely computer assisted data checks. In CAPI mode, interviewers record responses of a personal interview in an electronic questionnaire form on a phone or tablet. In CATI mode, Computer Assisted Telephone Interviewing, interviews are conducted over the phone and responses recorded in an electronic questionnaire on a computer, tablet or phone. "

# New row 13 - add the "Survey" concept (term only, definition left blank)
$ws.Range("A13").Value = "Survey"

# Reflect the author's final view/selection state
$ws.Range("A14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
